$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 33 so existing rows 33-39 shift down to 34-40,
# carrying their formatting (date style) along with them.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record.
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44642
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(34, 4).NumberFormat
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 100112030
$ws.Cells.Item(33, 7).Value = "Poroto granado"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 80
$ws.Cells.Item(33, 11).Value = 25000
$ws.Cells.Item(33, 12).Value = 27000
$ws.Cells.Item(33, 13).Value = 26000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Región Metropolitana"
$ws.Cells.Item(33, 16).Value = 1040
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
